$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.622.80'
$ws.Range("E2").Value = '  +0.40%  '
$ws.Range("D3").Value = '2.506.85'
$ws.Range("E3").Value = '  +0.02%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '575.11'
$ws.Range("E5").Value = '  +0.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.54'
$ws.Range("E6").Value = '  +0.42%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.513'
$ws.Range("E8").Value = '  -0.46%  '
$ws.Range("D9").Value = '2.505.00'
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +2.01%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.167'
$ws.Range("E11").Value = '  -0.37%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.357'
$ws.Range("E12").Value = '  +4.67%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.93'
$ws.Range("E13").Value = '  +2.05%  '
$ws.Range("D14").Value = '2.967.24'
$ws.Range("E14").Value = '  +0.09%  '
$ws.Range("D15").Value = '69.501.41'
$ws.Range("E15").Value = '  +0.30%  '
$ws.Range("E16").Value = '  +1.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '24.85'
$ws.Range("E17").Value = '  +0.62%  '
$ws.Range("D18").Value = '2.507.14'
$ws.Range("E18").Value = '  -0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.23'
$ws.Range("E19").Value = '  -1.37%  '
$ws.Range("E20").Value = '  -3.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '349.73'
$ws.Range("E21").Value = '  +0.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.92'
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.96'
$ws.Range("E23").Value = '  +0.50%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.99'
$ws.Range("E25").Value = '  +2.49%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.95'
$ws.Range("E26").Value = '  -0.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.82'
$ws.Range("E27").Value = '  -1.17%  '
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.34%  '
$ws.Range("D30").Value = '0.0₃0891'
$ws.Range("E30").Value = '  -0.62%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  +0.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '458.79'
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.73'
$ws.Range("E34").Value = '  -0.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.03'
$ws.Range("E36").Value = '  +4.27%  '
$ws.Range("E37").Value = '  +0.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.06'
$ws.Range("E38").Value = '  +0.74%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.47'
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("E40").Value = '  +0.00%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.319'
$ws.Range("E41").Value = '  +0.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.68'
$ws.Range("E42").Value = '  -0.91%  '
$ws.Range("E43").Value = '  +0.09%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '38.15'
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.21'
$ws.Range("E45").Value = '  -4.26%  '
$ws.Range("E46").Value = '  -7.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '142.07'
$ws.Range("E47").Value = '  -0.59%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.46'
$ws.Range("E48").Value = '  -1.36%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.519'
$ws.Range("E49").Value = '  -1.48%  '
$ws.Range("E50").Value = '  +0.20%  '
$ws.Range("E51").Value = '  -1.09%  '
